{"js": "// Replace each old math-fact text with its updated counterpart.\n// All 100 cell values in the table are unique, so an exact, case-sensitive,\n// whole-match search for each old string unambiguously finds the single run\n// that needs to change.\nconst replacements = [\n  [\"36+25=61\", \"47+24=71\"],\n  [\"36+6=42\", \"80-65=15\"],\n  [\"81-63=18\", \"20+30=50\"],\n  [\"75-56=19\", \"90+1=91\"],\n  [\"22+69=91\", \"34-7=27\"],\n  [\"94-46=48\", \"49+34=83\"],\n  [\"9+51=60\", \"15+17=32\"],\n  [\"67+9=76\", \"8+85=93\"],\n  [\"89+1=90\", \"84-6=78\"],\n  [\"99-70=29\", \"59-14=45\"],\n  [\"57-55=2\", \"24-9=15\"],\n  [\"39-11=28\", \"18+42=60\"],\n  [\"89-9=80\", \"53-19=34\"],\n  [\"50+12=62\", \"1+49=50\"],\n  [\"68-19=49\", \"65-3=62\"],\n  [\"94-49=45\", \"8+90=98\"],\n  [\"49-9=40\", \"36+27=63\"],\n  [\"44+31=75\", \"8+25=33\"],\n  [\"44-14=30\", \"0+78=78\"],\n  [\"80-56=24\", \"20+13=33\"],\n  [\"44+48=92\", \"30+15=45\"],\n  [\"97-71=26\", \"34-20=14\"],\n  [\"17+10=27\", \"28+16=44\"],\n  [\"85-82=3\", \"2+53=55\"],\n  [\"19+22=41\", \"1+20=21\"],\n  [\"95-3=92\", \"62-14=48\"],\n  [\"35-13=22\", \"29+18=47\"],\n  [\"60+14=74\", \"82-75=7\"],\n  [\"55-6=49\", \"93-78=15\"],\n  [\"63-55=8\", \"46-35=11\"],\n  [\"18-11=7\", \"45+34=79\"],\n  [\"33-11=22\", \"57+30=87\"],\n  [\"57+16=73\", \"86-78=8\"],\n  [\"28+64=92\", \"43-17=26\"],\n  [\"20+12=32\", \"53-13=40\"],\n  [\"4+36=40\", \"65-47=18\"],\n  [\"23+28=51\", \"30-21=9\"],\n  [\"89-37=52\", \"4+66=70\"],\n  [\"20+70=90\", \"10+34=44\"],\n  [\"72-52=20\", \"30+31=61\"],\n  [\"73-20=53\", \"19+63=82\"],\n  [\"39+14=53\", \"87-11=76\"],\n  [\"47-16=31\", \"9+1=10\"],\n  [\"11+86=97\", \"27+51=78\"],\n  [\"21+67=88\", \"82-3=79\"],\n  [\"54+11=65\", \"25-7=18\"],\n  [\"57-41=16\", \"9+24=33\"],\n  [\"37+26=63\", \"59-38=21\"],\n  [\"92-64=28\", \"17-12=5\"],\n  [\"49+46=95\", \"43+50=93\"],\n  [\"3+7=10\", \"31+64=95\"],\n  [\"83+2=85\", \"62+30=92\"],\n  [\"87-49=38\", \"51-30=21\"],\n  [\"66-33=33\", \"85-83=2\"],\n  [\"30-20=10\", \"92-84=8\"],\n  [\"55+39=94\", \"59-11=48\"],\n  [\"20+64=84\", \"33-8=25\"],\n  [\"31+29=60\", \"28-17=11\"],\n  [\"85-8=77\", \"12-11=1\"],\n  [\"95-29=66\", \"62-40=22\"],\n  [\"88-37=51\", \"36+53=89\"],\n  [\"47+48=95\", \"39+37=76\"],\n  [\"72+13=85\", \"1+76=77\"],\n  [\"21-20=1\", \"85-30=55\"],\n  [\"40-6=34\", \"10+32=42\"],\n  [\"95-13=82\", \"98-83=15\"],\n  [\"39-7=32\", \"61-38=23\"],\n  [\"72-12=60\", \"83-50=33\"],\n  [\"23+38=61\", \"48-42=6\"],\n  [\"48-12=36\", \"69+15=84\"],\n  [\"67-58=9\", \"49-29=20\"],\n  [\"48-27=21\", \"10+61=71\"],\n  [\"98-65=33\", \"74-14=60\"],\n  [\"88-25=63\", \"29-1=28\"],\n  [\"70-56=14\", \"30+23=53\"],\n  [\"56-49=7\", \"96-87=9\"],\n  [\"93-7=86\", \"4+41=45\"],\n  [\"48-32=16\", \"94-5=89\"],\n  [\"5+17=22\", \"76+6=82\"],\n  [\"61-45=16\", \"64+22=86\"],\n  [\"98-30=68\", \"99-29=70\"],\n  [\"58-9=49\", \"10+58=68\"],\n  [\"50+6=56\", \"37+4=41\"],\n  [\"82-4=78\", \"14+55=69\"],\n  [\"41-35=6\", \"70-16=54\"],\n  [\"83+11=94\", \"56-54=2\"],\n  [\"67-30=37\", \"86-45=41\"],\n  [\"91-85=6\", \"4+65=69\"],\n  [\"48-24=24\", \"83-5=78\"],\n  [\"3+88=91\", \"66-49=17\"],\n  [\"97-39=58\", \"94-71=23\"],\n  [\"57-44=13\", \"66+16=82\"],\n  [\"2+42=44\", \"8+56=64\"],\n  [\"55+1=56\", \"60-29=31\"],\n  [\"86-39=47\", \"7+62=69\"],\n  [\"30+57=87\", \"60+23=83\"],\n  [\"31+33=64\", \"14+80=94\"],\n  [\"72-22=50\", \"79-48=31\"],\n  [\"15+29=44\", \"75-58=17\"],\n  [\"89-23=66\", \"39+23=62\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load('items');\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text: ${oldText}`);\n  }\n\n  found.items[0].insertText(newText, 'Replace');\n}\n\nawait context.sync();", "ps1": "# Replace each old math-fact string with its updated counterpart.\n# All 100 cell values in the table are unique, so a case-sensitive\n# Find/Replace of each exact old string unambiguously hits the single run\n# that needs to change.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"36+25=61\", \"47+24=71\"),\n    @(\"36+6=42\", \"80-65=15\"),\n    @(\"81-63=18\", \"20+30=50\"),\n    @(\"75-56=19\", \"90+1=91\"),\n    @(\"22+69=91\", \"34-7=27\"),\n    @(\"94-46=48\", \"49+34=83\"),\n    @(\"9+51=60\", \"15+17=32\"),\n    @(\"67+9=76\", \"8+85=93\"),\n    @(\"89+1=90\", \"84-6=78\"),\n    @(\"99-70=29\", \"59-14=45\"),\n    @(\"57-55=2\", \"24-9=15\"),\n    @(\"39-11=28\", \"18+42=60\"),\n    @(\"89-9=80\", \"53-19=34\"),\n    @(\"50+12=62\", \"1+49=50\"),\n    @(\"68-19=49\", \"65-3=62\"),\n    @(\"94-49=45\", \"8+90=98\"),\n    @(\"49-9=40\", \"36+27=63\"),\n    @(\"44+31=75\", \"8+25=33\"),\n    @(\"44-14=30\", \"0+78=78\"),\n    @(\"80-56=24\", \"20+13=33\"),\n    @(\"44+48=92\", \"30+15=45\"),\n    @(\"97-71=26\", \"34-20=14\"),\n    @(\"17+10=27\", \"28+16=44\"),\n    @(\"85-82=3\", \"2+53=55\"),\n    @(\"19+22=41\", \"1+20=21\"),\n    @(\"95-3=92\", \"62-14=48\"),\n    @(\"35-13=22\", \"29+18=47\"),\n    @(\"60+14=74\", \"82-75=7\"),\n    @(\"55-6=49\", \"93-78=15\"),\n    @(\"63-55=8\", \"46-35=11\"),\n    @(\"18-11=7\", \"45+34=79\"),\n    @(\"33-11=22\", \"57+30=87\"),\n    @(\"57+16=73\", \"86-78=8\"),\n    @(\"28+64=92\", \"43-17=26\"),\n    @(\"20+12=32\", \"53-13=40\"),\n    @(\"4+36=40\", \"65-47=18\"),\n    @(\"23+28=51\", \"30-21=9\"),\n    @(\"89-37=52\", \"4+66=70\"),\n    @(\"20+70=90\", \"10+34=44\"),\n    @(\"72-52=20\", \"30+31=61\"),\n    @(\"73-20=53\", \"19+63=82\"),\n    @(\"39+14=53\", \"87-11=76\"),\n    @(\"47-16=31\", \"9+1=10\"),\n    @(\"11+86=97\", \"27+51=78\"),\n    @(\"21+67=88\", \"82-3=79\"),\n    @(\"54+11=65\", \"25-7=18\"),\n    @(\"57-41=16\", \"9+24=33\"),\n    @(\"37+26=63\", \"59-38=21\"),\n    @(\"92-64=28\", \"17-12=5\"),\n    @(\"49+46=95\", \"43+50=93\"),\n    @(\"3+7=10\", \"31+64=95\"),\n    @(\"83+2=85\", \"62+30=92\"),\n    @(\"87-49=38\", \"51-30=21\"),\n    @(\"66-33=33\", \"85-83=2\"),\n    @(\"30-20=10\", \"92-84=8\"),\n    @(\"55+39=94\", \"59-11=48\"),\n    @(\"20+64=84\", \"33-8=25\"),\n    @(\"31+29=60\", \"28-17=11\"),\n    @(\"85-8=77\", \"12-11=1\"),\n    @(\"95-29=66\", \"62-40=22\"),\n    @(\"88-37=51\", \"36+53=89\"),\n    @(\"47+48=95\", \"39+37=76\"),\n    @(\"72+13=85\", \"1+76=77\"),\n    @(\"21-20=1\", \"85-30=55\"),\n    @(\"40-6=34\", \"10+32=42\"),\n    @(\"95-13=82\", \"98-83=15\"),\n    @(\"39-7=32\", \"61-38=23\"),\n    @(\"72-12=60\", \"83-50=33\"),\n    @(\"23+38=61\", \"48-42=6\"),\n    @(\"48-12=36\", \"69+15=84\"),\n    @(\"67-58=9\", \"49-29=20\"),\n    @(\"48-27=21\", \"10+61=71\"),\n    @(\"98-65=33\", \"74-14=60\"),\n    @(\"88-25=63\", \"29-1=28\"),\n    @(\"70-56=14\", \"30+23=53\"),\n    @(\"56-49=7\", \"96-87=9\"),\n    @(\"93-7=86\", \"4+41=45\"),\n    @(\"48-32=16\", \"94-5=89\"),\n    @(\"5+17=22\", \"76+6=82\"),\n    @(\"61-45=16\", \"64+22=86\"),\n    @(\"98-30=68\", \"99-29=70\"),\n    @(\"58-9=49\", \"10+58=68\"),\n    @(\"50+6=56\", \"37+4=41\"),\n    @(\"82-4=78\", \"14+55=69\"),\n    @(\"41-35=6\", \"70-16=54\"),\n    @(\"83+11=94\", \"56-54=2\"),\n    @(\"67-30=37\", \"86-45=41\"),\n    @(\"91-85=6\", \"4+65=69\"),\n    @(\"48-24=24\", \"83-5=78\"),\n    @(\"3+88=91\", \"66-49=17\"),\n    @(\"97-39=58\", \"94-71=23\"),\n    @(\"57-44=13\", \"66+16=82\"),\n    @(\"2+42=44\", \"8+56=64\"),\n    @(\"55+1=56\", \"60-29=31\"),\n    @(\"86-39=47\", \"7+62=69\"),\n    @(\"30+57=87\", \"60+23=83\"),\n    @(\"31+33=64\", \"14+80=94\"),\n    @(\"72-22=50\", \"79-48=31\"),\n    @(\"15+29=44\", \"75-58=17\"),\n    @(\"89-23=66\", \"39+23=62\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute(\n        $oldText,  # FindText\n        $true,     # MatchCase\n        $false,    # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap (wdFindContinue)\n        $false,    # Format\n        $newText,  # ReplaceWith\n        2          # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"Could not find text: $oldText\"\n    }\n}\n"}
